$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.996.35'
$ws.Range('E2').Value = '  +1.86%  '
$ws.Range('D3').Value = '1.670.68'
$ws.Range('E3').Value = '  +2.75%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'215.93"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').Value = "'0.512"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.94%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +1.96%  '
$ws.Range('E9').Value = '  +1.10%  '
$ws.Range('E10').Value = '  +4.53%  '
$ws.Range('D11').Value = "'0.0892"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.48%  '
$ws.Range('D12').Value = '1.908.65'
$ws.Range('E12').Value = '  +2.90%  '
$ws.Range('D13').Value = '1.670.54'
$ws.Range('E13').Value = '  +2.75%  '
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('D15').Value = "'65.74"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.70%  '
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('D17').Value = '27.025.75'
$ws.Range('E17').Value = '  +1.96%  '
$ws.Range('D18').Value = "'235.13"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.06%  '
$ws.Range('E19').Value = '  +1.21%  '
$ws.Range('D20').Value = "'7.71"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.70%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('E22').Value = '  +3.29%  '
$ws.Range('E23').Value = '  +1.47%  '
$ws.Range('E24').Value = '  +1.11%  '
$ws.Range('D25').Value = "'145.42"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('E26').Value = '  +1.14%  '
$ws.Range('E27').Value = '  +0.59%  '
$ws.Range('D28').Value = "'15.86"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  +1.32%  '
$ws.Range('D32').Value = "'3.32"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.03%  '
$ws.Range('D33').Value = '1.452.57'
$ws.Range('E33').Value = '  -4.37%  '
$ws.Range('E34').Value = '  +5.23%  '
$ws.Range('D35').Value = "'1.60"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.22%  '
$ws.Range('D36').Value = "'2.42"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = "'0.568"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').Value = "'0.890"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.57%  '
$ws.Range('D39').Value = "'0.0169"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.74%  '
$ws.Range('E40').Value = '  +3.54%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('E42').Value = '  +11.33%  '
$ws.Range('D43').Value = "'2.28"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.93%  '
$ws.Range('D44').Value = "'65.74"
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Value = '1.816.26'
$ws.Range('E45').Value = '  +2.84%  '
$ws.Range('E46').Value = '  +2.57%  '
$ws.Range('D47').Value = "'90.24"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('E48').Value = '  +1.25%  '
$ws.Range('E49').Value = '  +3.99%  '
$ws.Range('E50').Value = '  +1.33%  '
$ws.Range('E51').Value = '  +1.44%  '
